$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GATED")
$ws1.Range("E5").Value = 85.7013
$ws1.Range("F5").Value = 85.9425
